$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H119").Value = 2297
$ws.Range("J119").Value = 2297
$ws.Range("L119").Value = 6891
$ws.Range("N119").Value = -16567

$ws.Range("H132").Value = 3748
$ws.Range("I132").Value = 4067.2068
$ws.Range("K132").Value = 12201.6204
$ws.Range("M132").Value = -9671.6204

$ws.Range("H135").Value = 660.7857
$ws.Range("I135").Value = 518.375
$ws.Range("J135").Value = 850.6667
$ws.Range("K135").Value = 4665.375
$ws.Range("L135").Value = 7656.0003
$ws.Range("M135").Value = -2130.375
$ws.Range("N135").Value = -12726.0003

$ws.Range("H138").Value = 13892775
$ws.Range("J138").Value = 24396410
$ws.Range("L138").Value = 73189230
$ws.Range("N138").Value = -73199510

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4442.533
$ws.Range("I2").Value = 4187
$ws.Range("K2").Value = 4187
$ws.Range("M2").Value = -4074

$ws.Range("H45").Value = 4077.182
$ws.Range("I45").Value = 3434.7144
$ws.Range("K45").Value = 3434.7144
$ws.Range("M45").Value = -3057.7144

$ws.Range("H61").Value = 3256.2068
$ws.Range("I61").Value = 2917.818
$ws.Range("K61").Value = 2917.818
$ws.Range("M61").Value = -2705.818

$ws.Range("H110").Value = 39909.223
$ws.Range("I110").Value = 39909.223
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 39909.223
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -37864.223
$ws.Range("N110").Value = ""

$ws.Range("H116").Value = 4442.533
$ws.Range("I116").Value = 4187
$ws.Range("K116").Value = 4187
$ws.Range("M116").Value = -1893

$ws.Range("H132").Value = 82153.13
$ws.Range("I132").Value = 5106.115
$ws.Range("J132").Value = 482797.6
$ws.Range("K132").Value = 15318.345
$ws.Range("L132").Value = 1448392.8
$ws.Range("M132").Value = -12788.345
$ws.Range("N132").Value = -1453452.8

$ws.Range("H136").Value = 3256.2068
$ws.Range("I136").Value = 2917.818
$ws.Range("K136").Value = 8753.454000000002
$ws.Range("M136").Value = -6203.454000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4442.533
$ws.Range("I3").Value = 4187
$ws.Range("K3").Value = 4187
$ws.Range("M3").Value = -4073

$ws.Range("H86").Value = 40724.875
$ws.Range("I86").Value = 27126.5
$ws.Range("J86").Value = 54323.25
$ws.Range("K86").Value = 27126.5
$ws.Range("L86").Value = 54323.25
$ws.Range("M86").Value = -26003.5
$ws.Range("N86").Value = -56569.25

$ws.Range("H89").Value = 40724.875
$ws.Range("I89").Value = 27126.5
$ws.Range("J89").Value = 54323.25
$ws.Range("K89").Value = 135632.5
$ws.Range("L89").Value = 271616.25
$ws.Range("M89").Value = -130016.5
$ws.Range("N89").Value = -282848.25

$ws.Range("H94").Value = 1902.4
$ws.Range("I94").Value = 1734.3334
$ws.Range("J94").Value = 2154.5
$ws.Range("K94").Value = 1734.3334
$ws.Range("L94").Value = 2154.5
$ws.Range("M94").Value = -1283.3334
$ws.Range("N94").Value = -3056.5

$ws.Range("H107").Value = 2684.6667
$ws.Range("I107").Value = 2389.0667
$ws.Range("K107").Value = 2389.0667
$ws.Range("M107").Value = -469.0666999999999

$ws.Range("H134").Value = 1145.9
$ws.Range("I134").Value = 1053.2222
$ws.Range("K134").Value = 3159.6666
$ws.Range("M134").Value = -624.6665999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 60000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = ""

$ws.Range("H31").Value = 2584.0981
$ws.Range("I31").Value = 1062.6578
$ws.Range("K31").Value = 1062.6578
$ws.Range("M31").Value = -767.6578

$ws.Range("H34").Value = 2584.0981
$ws.Range("I34").Value = 1062.6578
$ws.Range("K34").Value = 1062.6578
$ws.Range("M34").Value = -860.6578

$ws.Range("H48").Value = 41199
$ws.Range("J48").Value = 41199
$ws.Range("L48").Value = 41199
$ws.Range("N48").Value = -42151

$ws.Range("H64").Value = 35000
$ws.Range("J64").Value = 35000
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496

$ws.Range("H67").Value = 35000
$ws.Range("J67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""

$ws.Range("H132").Value = 3827.2632
$ws.Range("I132").Value = 3575.3635
$ws.Range("J132").Value = 5489.8
$ws.Range("K132").Value = 10726.0905
$ws.Range("L132").Value = 16469.4
$ws.Range("M132").Value = -8196.0905
$ws.Range("N132").Value = -21529.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 800
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""

$ws.Range("H54").Value = 4990
$ws.Range("J54").Value = 4990
$ws.Range("L54").Value = 14970
$ws.Range("N54").Value = -16088

$ws.Range("H97").Value = 797.7143
$ws.Range("I97").Value = 1091
$ws.Range("J97").Value = 406.66666
$ws.Range("K97").Value = 3273
$ws.Range("L97").Value = 1219.99998
$ws.Range("M97").Value = -2777
$ws.Range("N97").Value = -2211.99998

$ws.Range("H109").Value = 8595.833000000001
$ws.Range("I109").Value = 21013
$ws.Range("J109").Value = 2387.25
$ws.Range("K109").Value = 63039
$ws.Range("L109").Value = 7161.75
$ws.Range("M109").Value = -61999
$ws.Range("N109").Value = -9241.75

$ws.Range("H131").Value = 1591.1818
$ws.Range("I131").Value = 1325.5
$ws.Range("J131").Value = 1743
$ws.Range("K131").Value = 3976.5
$ws.Range("L131").Value = 5229
$ws.Range("M131").Value = 1063.5
$ws.Range("N131").Value = -15309

$ws.Range("H132").Value = 3850.5
$ws.Range("I132").Value = 1449
$ws.Range("K132").Value = 13041
$ws.Range("M132").Value = -10511

$ws.Range("H133").Value = 3033
$ws.Range("J133").Value = 3033
$ws.Range("L133").Value = 9099
$ws.Range("N133").Value = -19219

$ws.Range("H134").Value = 1369.9474
$ws.Range("I134").Value = 1174.2941
$ws.Range("J134").Value = 3033
$ws.Range("K134").Value = 3522.8823
$ws.Range("L134").Value = 9099
$ws.Range("M134").Value = 1547.1177
$ws.Range("N134").Value = -19239

$ws.Range("H137").Value = 1696.5
$ws.Range("J137").Value = 2365
$ws.Range("L137").Value = 7095
$ws.Range("N137").Value = -17295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""

$ws.Range("H33").Value = 1012999.3
$ws.Range("I33").Value = 19499
$ws.Range("K33").Value = 19499
$ws.Range("M33").Value = -19247

$ws.Range("H97").Value = 1710.8334
$ws.Range("I97").Value = 1172.5714
$ws.Range("J97").Value = 2464.4
$ws.Range("K97").Value = 1172.5714
$ws.Range("L97").Value = 2464.4
$ws.Range("M97").Value = -676.5714
$ws.Range("N97").Value = -3456.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 41497.5

$ws.Range("H61").Value = 22729956
$ws.Range("I61").Value = 33335814
$ws.Range("K61").Value = 33335814
$ws.Range("M61").Value = -33335612

$ws.Range("H113").Value = 22729956
$ws.Range("I113").Value = 33335814
$ws.Range("K113").Value = 33335814
$ws.Range("M113").Value = -33333644

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 80
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 80
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = -360

$ws.Range("H14").Value = 6545.423
$ws.Range("I14").Value = 131.22223
$ws.Range("K14").Value = 131.22223
$ws.Range("M14").Value = 36.77777

$ws.Range("H15").Value = 9200
$ws.Range("I15").Value = 5000
$ws.Range("K15").Value = 5000
$ws.Range("M15").Value = -4712

$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = ""

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""

$ws.Range("H52").Value = 36155.844
$ws.Range("J52").Value = 38747.562
$ws.Range("L52").Value = 38747.562
$ws.Range("N52").Value = -39199.562

$ws.Range("H70").Value = 31000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""

$ws.Range("H73").Value = 31000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""

$ws.Range("H96").Value = 3561.1304
$ws.Range("I96").Value = 3186.1667
$ws.Range("K96").Value = 3186.1667
$ws.Range("M96").Value = -1813.1667
